$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused helper notes in column E
$ws.Range("E31").ClearContents()
$ws.Range("E41").ClearContents()

# Insert a new blank row before the old "ukupno" row, pushing the totals down
$ws.Rows.Item(54).Insert()

# Make the summary labels bold
$ws.Range("A52").Font.Bold = $true
$ws.Range("A53").Font.Bold = $true

# Update the active selection/view
[void]$ws.Range("E21").Select()

Write-Output "done"
